$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "88-70=18"  # was: 71+10=81
$t.Cell(1, 2).Range.Text = "5+81=86"  # was: 73-37=36
$t.Cell(1, 3).Range.Text = "9+10=19"  # was: 2+36=38
$t.Cell(1, 4).Range.Text = "88+8=96"  # was: 75-14=61
$t.Cell(1, 5).Range.Text = "67-65=2"  # was: 22-2=20
$t.Cell(2, 1).Range.Text = "17+75=92"  # was: 28-13=15
$t.Cell(2, 2).Range.Text = "77-47=30"  # was: 34-4=30
$t.Cell(2, 3).Range.Text = "3+88=91"  # was: 47-32=15
$t.Cell(2, 4).Range.Text = "95-66=29"  # was: 99-69=30
$t.Cell(2, 5).Range.Text = "23+4=27"  # was: 80+18=98
$t.Cell(3, 1).Range.Text = "8+80=88"  # was: 47-31=16
$t.Cell(3, 2).Range.Text = "38+18=56"  # was: 80+13=93
$t.Cell(3, 3).Range.Text = "54-27=27"  # was: 74-14=60
$t.Cell(3, 4).Range.Text = "92-30=62"  # was: 15-4=11
$t.Cell(3, 5).Range.Text = "76-42=34"  # was: 63-1=62
$t.Cell(4, 1).Range.Text = "2+78=80"  # was: 47+20=67
$t.Cell(4, 2).Range.Text = "31+68=99"  # was: 99+0=99
$t.Cell(4, 3).Range.Text = "53-11=42"  # was: 68+20=88
$t.Cell(4, 4).Range.Text = "11+54=65"  # was: 51-25=26
$t.Cell(4, 5).Range.Text = "1+80=81"  # was: 4+66=70
$t.Cell(5, 1).Range.Text = "69+0=69"  # was: 21+33=54
$t.Cell(5, 2).Range.Text = "76-37=39"  # was: 98-85=13
$t.Cell(5, 3).Range.Text = "12+68=80"  # was: 12+51=63
$t.Cell(5, 4).Range.Text = "41+31=72"  # was: 14+13=27
$t.Cell(5, 5).Range.Text = "62-21=41"  # was: 65+2=67
$t.Cell(6, 1).Range.Text = "39+3=42"  # was: 80-2=78
$t.Cell(6, 2).Range.Text = "21+67=88"  # was: 59-0=59
$t.Cell(6, 3).Range.Text = "31+37=68"  # was: 28+49=77
$t.Cell(6, 4).Range.Text = "12+4=16"  # was: 48+28=76
$t.Cell(6, 5).Range.Text = "47+32=79"  # was: 16+58=74
$t.Cell(7, 1).Range.Text = "38+8=46"  # was: 13+67=80
$t.Cell(7, 2).Range.Text = "25-11=14"  # was: 55-51=4
$t.Cell(7, 3).Range.Text = "47+41=88"  # was: 97-50=47
$t.Cell(7, 4).Range.Text = "91-41=50"  # was: 54+41=95
$t.Cell(7, 5).Range.Text = "96-8=88"  # was: 43-5=38
$t.Cell(8, 1).Range.Text = "33+54=87"  # was: 43+44=87
$t.Cell(8, 2).Range.Text = "82-58=24"  # was: 56+42=98
$t.Cell(8, 3).Range.Text = "49+11=60"  # was: 5+57=62
$t.Cell(8, 4).Range.Text = "10+86=96"  # was: 57+40=97
$t.Cell(8, 5).Range.Text = "71+19=90"  # was: 83-4=79
$t.Cell(9, 1).Range.Text = "45+32=77"  # was: 97-41=56
$t.Cell(9, 2).Range.Text = "10+53=63"  # was: 34+2=36
$t.Cell(9, 3).Range.Text = "81-16=65"  # was: 33+12=45
$t.Cell(9, 4).Range.Text = "98-41=57"  # was: 7+55=62
$t.Cell(9, 5).Range.Text = "24+1=25"  # was: 4-1=3
$t.Cell(10, 1).Range.Text = "6+85=91"  # was: 29+29=58
$t.Cell(10, 2).Range.Text = "68-11=57"  # was: 0+99=99
$t.Cell(10, 3).Range.Text = "48-19=29"  # was: 9+84=93
$t.Cell(10, 4).Range.Text = "89-46=43"  # was: 45+36=81
$t.Cell(10, 5).Range.Text = "66-32=34"  # was: 91-57=34
$t.Cell(11, 1).Range.Text = "69-1=68"  # was: 83+6=89
$t.Cell(11, 2).Range.Text = "12+87=99"  # was: 17+17=34
$t.Cell(11, 3).Range.Text = "60+27=87"  # was: 47+46=93
$t.Cell(11, 4).Range.Text = "28+71=99"  # was: 62-56=6
$t.Cell(11, 5).Range.Text = "87-74=13"  # was: 58-29=29
$t.Cell(12, 1).Range.Text = "61-6=55"  # was: 12+60=72
$t.Cell(12, 2).Range.Text = "61+3=64"  # was: 94-9=85
$t.Cell(12, 3).Range.Text = "81-45=36"  # was: 92-70=22
$t.Cell(12, 4).Range.Text = "49-35=14"  # was: 85-80=5
$t.Cell(12, 5).Range.Text = "74+5=79"  # was: 50+5=55
$t.Cell(13, 1).Range.Text = "26-7=19"  # was: 27-4=23
$t.Cell(13, 2).Range.Text = "39+46=85"  # was: 55+26=81
$t.Cell(13, 3).Range.Text = "3+44=47"  # was: 83-83=0
$t.Cell(13, 4).Range.Text = "56-0=56"  # was: 69-45=24
$t.Cell(13, 5).Range.Text = "1+13=14"  # was: 75-65=10
$t.Cell(14, 1).Range.Text = "67-60=7"  # was: 78-74=4
$t.Cell(14, 2).Range.Text = "80-4=76"  # was: 44+18=62
$t.Cell(14, 3).Range.Text = "63-41=22"  # was: 42+24=66
$t.Cell(14, 4).Range.Text = "70-32=38"  # was: 48-21=27
$t.Cell(14, 5).Range.Text = "82+11=93"  # was: 77-41=36
$t.Cell(15, 1).Range.Text = "95-38=57"  # was: 33+25=58
$t.Cell(15, 2).Range.Text = "77+10=87"  # was: 68+30=98
$t.Cell(15, 3).Range.Text = "57-0=57"  # was: 52+40=92
$t.Cell(15, 4).Range.Text = "6+6=12"  # was: 54+30=84
$t.Cell(15, 5).Range.Text = "4+76=80"  # was: 29+47=76
$t.Cell(16, 1).Range.Text = "1+87=88"  # was: 13+24=37
$t.Cell(16, 2).Range.Text = "85+11=96"  # was: 78+4=82
$t.Cell(16, 3).Range.Text = "23+5=28"  # was: 23+76=99
$t.Cell(16, 4).Range.Text = "87-75=12"  # was: 93-87=6
$t.Cell(16, 5).Range.Text = "55+14=69"  # was: 38+14=52
$t.Cell(17, 1).Range.Text = "37+17=54"  # was: 42-4=38
$t.Cell(17, 2).Range.Text = "96-5=91"  # was: 23+30=53
$t.Cell(17, 3).Range.Text = "21+45=66"  # was: 87+7=94
$t.Cell(17, 4).Range.Text = "23+28=51"  # was: 6+50=56
$t.Cell(17, 5).Range.Text = "95-75=20"  # was: 10+10=20
$t.Cell(18, 1).Range.Text = "10+44=54"  # was: 71+19=90
$t.Cell(18, 2).Range.Text = "88-47=41"  # was: 13+77=90
$t.Cell(18, 3).Range.Text = "17+32=49"  # was: 39+8=47
$t.Cell(18, 4).Range.Text = "49-33=16"  # was: 36-11=25
$t.Cell(18, 5).Range.Text = "42-21=21"  # was: 28-27=1
$t.Cell(19, 1).Range.Text = "68-41=27"  # was: 85-30=55
$t.Cell(19, 2).Range.Text = "57-54=3"  # was: 33+41=74
$t.Cell(19, 3).Range.Text = "87-34=53"  # was: 54-50=4
$t.Cell(19, 4).Range.Text = "7+60=67"  # was: 48-9=39
$t.Cell(19, 5).Range.Text = "36+12=48"  # was: 74+24=98
$t.Cell(20, 1).Range.Text = "87+1=88"  # was: 77-38=39
$t.Cell(20, 2).Range.Text = "48+35=83"  # was: 64-63=1
$t.Cell(20, 3).Range.Text = "47+49=96"  # was: 97-52=45
$t.Cell(20, 4).Range.Text = "96-92=4"  # was: 88-11=77
$t.Cell(20, 5).Range.Text = "58+37=95"  # was: 21+30=51
